$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert the new "PO Ref / Contractor" row at row 2 --------------
# Row 1 has the exact format pattern we need (plain / s1 / plain / s1), so
# clone its formatting down onto row 2 before writing the new values.
$ws.Range("A1:E1").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)
$ws.Range("C2").Clear()
$ws.Range("A2").Value = "PO Ref"
$ws.Range("B2").Value = "GVT000ECI24000829"
$ws.Range("D2").Value = "Contractor"
$ws.Range("E2").Value = "PALO IT"

# --- Step 2: row 3 becomes the "PO Date" row (was "PO Ref") -----------------
$ws.Range("A3").Value = "PO Date"
$ws.Range("B3").Value = "1 May 24 - 30 Apr 25"
$ws.Range("D3").Clear()
$ws.Range("E3").Clear()

# --- Step 3: the old standalone "PO Date" row (row 4) is now empty ----------
$ws.Range("A4:B4").Clear()

# --- Step 4: row 7 keeps only Role Specialization; Group/Specialization ----
# moves onto its own row 8.
$ws.Range("D7").Clear()
$ws.Range("E7").ClearContents()

$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B6").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("E6").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Group/Specialization"
$ws.Range("B8").Value = "Consultant"
$ws.Range("E8").ClearContents()

# --- Step 5: 2025-01-02 (row 11) flips from Annual Leave to At Work --------
$ws.Range("A12:H12").Copy()
$ws.Range("A11:H11").PasteSpecial(-4122)
$ws.Range("C11").Value = 1
$ws.Range("G11").Value = 0
$ws.Range("H11").ClearContents()

# --- Step 6: 2025-01-13 (row 22) flips from At Work to Sick Leave ----------
$ws.Range("A13:H13").Copy()
$ws.Range("A22:H22").PasteSpecial(-4122)
$ws.Range("C22").Value = 0
$ws.Range("E22").Value = 1
$ws.Range("H22").Value = "Sick Leave"

# --- Step 7: update the Totals row to reflect the two swapped days --------
$ws.Range("E41").Value = 6
$ws.Range("G41").Value = 0
